# Update countries & provincias Spain
# - Reorder Honduras/Taiwan in the country list (Honduras' updated data now
#   ranks just ahead of Taiwan in the table, row 97/98 swap countries+data).
# - Refresh the "Datos actualizados" timestamp.
# - Refresh the case counters for the countries/regions whose figures changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last refreshed" timestamp (row 1) -----------------------------
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 11 de Abril de 2020 a las 09:22"

# --- Helper to push a full data row (B..H) ---------------------------------
function Set-RowValues {
    param($row, $b, $c, $d, $e, $f, $g, $h)
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# Row 4 - Estados Unidos
Set-RowValues 4 503177 301 27314 457102 10917 14 18761

# Row 19 - Austria
Set-RowValues 19 13586 26 6064 7203 261 0 319

# Row 34 - Chequia
Set-RowValues 34 5735 3 370 5242 92 4 123

# Row 73 - Armenia
Set-RowValues 73 966 29 173 780 30 1 13

# Row 75 - Kazajistan
Set-RowValues 75 859 47 64 785 21 0 10

# Rows 97/98 - Honduras moves ahead of Taiwan with its refreshed figures,
# Taiwan keeps its previous totals one row below.
$ws.Cells.Item(97, 1).Value = "Honduras"
Set-RowValues 97 392 10 7 361 10 1 24

$ws.Cells.Item(98, 1).Value = "Taiwan"
Set-RowValues 98 385 3 99 280 0 0 6

# Row 110 - Montenegro
Set-RowValues 110 262 5 4 256 7 0 2
